$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.337.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.65%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.976.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.65%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "502.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.37"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +8.47%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("E8").Value = "  +6.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +11.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.106"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +10.92%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.350"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.78%  "

$ws.Range("E12").Value = "  +3.00%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.482.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +11.37%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "56.320.95"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.58%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000150"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +13.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.974.81"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.65"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +11.04%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.32"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.77%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.35%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +5.85%  "

$ws.Range("E22").Value = "  -0.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.470"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +6.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.997"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("E26").Value = "  +5.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0₃0886"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.12%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +10.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +12.56%  "

$ws.Range("E31").Value = "  +7.06%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +10.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.42"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +7.88%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "156.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +8.25%  "

$ws.Range("E35").Value = "  +6.70%  "

$ws.Range("E36").Value = "  +4.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.53"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0672"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "22.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.00%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.007.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.74%  "

$ws.Range("E41").Value = "  +0.04%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.42%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.641"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.33%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.245.73"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +10.09%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.984"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.60%  "

$ws.Range("E46").Value = "  +7.37%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.76%  "

$ws.Range("E48").Value = "  +25.65%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0234"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.61%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.88"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.84%  "
